# Insert a new row at position 72, shifting existing rows 72:180 down to 73:181
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("72:72").Insert()

# Populate the newly inserted row 72 with the new record's data.
# Columns A,B,C,E,F,G,H,I,R are identical metadata for every row in this sheet,
# so copy them from the row immediately below (row 73), which holds the data
# that used to live in row 72 before the insert.
$ws.Range("A72").Value = $ws.Range("A73").Value()
$ws.Range("B72").Value = $ws.Range("B73").Value()
$ws.Range("C72").Value = $ws.Range("C73").Value()
$ws.Range("D72").Value = 44580
$ws.Range("E72").Value = $ws.Range("E73").Value()
$ws.Range("F72").Value = $ws.Range("F73").Value()
$ws.Range("G72").Value = $ws.Range("G73").Value()
$ws.Range("H72").Value = $ws.Range("H73").Value()
$ws.Range("I72").Value = $ws.Range("I73").Value()
$ws.Range("J72").Value = 35
$ws.Range("K72").Value = 16000
$ws.Range("L72").Value = 16000
$ws.Range("M72").Value = 16000
$ws.Range("N72").Value = "$/caja 50 unidades"
$ws.Range("O72").Value = "Región de O'Higgins"
$ws.Range("P72").Value = 320
$ws.Range("Q72").Value = 50
$ws.Range("R72").Value = $ws.Range("R73").Value()
